$wb = $excel.ActiveWorkbook

# Sheet "Gesamtergebnis" now gets the values/formatting for the Total row's
# balance columns (C3/D3), and becomes the active/selected sheet
# (instead of "Tagesergebnisse").
$wsGes = $wb.Worksheets.Item("Gesamtergebnis")

# Set numeric totals instead of "N/A" text for balance columns
$wsGes.Range("C3").Value = 0
$wsGes.Range("D3").Value = 0

# Update row height for row 3
$wsGes.Rows.Item(3).RowHeight = 13.8

# Update selection on Gesamtergebnis sheet
$wsGes.Range("C7").Select()

# Make Gesamtergebnis the active sheet / tab
$wsGes.Activate()
